# Apply the "deploy to gh-pages" metadata refresh to the FHIR
# StructureDefinition workbook:
#   - URL host ibm.com -> linuxforhealth.org
#   - Version 7.0.0 -> 8.0.0
#   - Date bump
#   - Publisher Alvearie Team -> LinuxForHealth Team
#   - Clear the stray Constraint(s) text on the root Extension row of
#     the Elements sheet (it is now only carried on the Extension.extension row)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-derived-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
